# Estonia Meistriliiga base update (12-06-2024 23:38)
#
# Two match rows had their data entered against the wrong fixture on the
# source feed. Fix it by swapping the full set of match data (id, teams,
# score, odds, ...) between the two affected row pairs, while leaving the
# running index in column A untouched (it always mirrors the sheet row).
#
#   - rows 4 and 5   (fixtures played on 2023-06-13, serial 45084.5)
#   - rows 169 and 170 (fixtures played on 2024-05-27, serial 45434.541...)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param(
        $Worksheet,
        [int]$RowA,
        [int]$RowB
    )

    # Column A (the running id) stays put - only B:AD (id, Div, Date,
    # HomeTeam, AwayTeam, scores, odds, ...) swap between the two rows.
    $rangeA = $Worksheet.Range("B$RowA" + ":AD$RowA")
    $rangeB = $Worksheet.Range("B$RowB" + ":AD$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-MatchRows -Worksheet $ws -RowA 4 -RowB 5
Swap-MatchRows -Worksheet $ws -RowA 169 -RowB 170
